$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply each cell update using an apostrophe-prefixed value so Excel
# stores the result as literal text (matching the original inlineStr
# cells) instead of auto-converting number-looking strings into
# numeric values, then reset the style so no stray quote-prefix /
# number-format style survives on the cell.
function Set-TextValue {
    param($ws, $addr, $value)
    $ws.Range($addr).Value = "'" + $value
    $ws.Range($addr).Style = "Normal"
}

Set-TextValue $ws "D2" "28.320.95"
Set-TextValue $ws "E2" "  -0.77%  "
Set-TextValue $ws "D3" "1.566.56"
Set-TextValue $ws "E3" "  +0.08%  "
Set-TextValue $ws "E4" "  -0.02%  "
Set-TextValue $ws "D5" "210.92"
Set-TextValue $ws "E5" "  -0.39%  "
Set-TextValue $ws "E6" "  -0.71%  "
Set-TextValue $ws "E7" "  +0.02%  "
Set-TextValue $ws "D8" "44.36"
Set-TextValue $ws "E8" "  -4.30%  "
Set-TextValue $ws "D9" "23.73"
Set-TextValue $ws "E9" "  -1.73%  "
Set-TextValue $ws "D10" "0.244"
Set-TextValue $ws "E10" "  -1.11%  "
Set-TextValue $ws "E11" "  -0.87%  "
Set-TextValue $ws "D12" "0.0893"
Set-TextValue $ws "E12" "  +1.13%  "
Set-TextValue $ws "D13" "1.790.07"
Set-TextValue $ws "E13" "  +0.06%  "
Set-TextValue $ws "D14" "1.557.72"
Set-TextValue $ws "E14" "  -0.50%  "
Set-TextValue $ws "D15" "3.66"
Set-TextValue $ws "E15" "  -0.41%  "
Set-TextValue $ws "D16" "28.322.91"
Set-TextValue $ws "E16" "  -0.75%  "
Set-TextValue $ws "E17" "  -1.23%  "
Set-TextValue $ws "D18" "61.06"
Set-TextValue $ws "E18" "  -1.45%  "
Set-TextValue $ws "D19" "227.55"
Set-TextValue $ws "E19" "  +0.20%  "
Set-TextValue $ws "D20" "7.38"
Set-TextValue $ws "E20" "  +0.97%  "
Set-TextValue $ws "E21" "  -2.10%  "
Set-TextValue $ws "E22" "  -0.01%  "
Set-TextValue $ws "D23" "3.93"
Set-TextValue $ws "E23" "  +1.97%  "
Set-TextValue $ws "D24" "8.93"
Set-TextValue $ws "E24" "  -2.23%  "
Set-TextValue $ws "E25" "  -0.59%  "
Set-TextValue $ws "D26" "150.58"
Set-TextValue $ws "E26" "  -0.21%  "
Set-TextValue $ws "D27" "14.89"
Set-TextValue $ws "E27" "  -0.45%  "
Set-TextValue $ws "E28" "  -0.21%  "
Set-TextValue $ws "E29" "  -1.63%  "
Set-TextValue $ws "E30" "  -0.02%  "
Set-TextValue $ws "E31" "  +3.13%  "
Set-TextValue $ws "E32" "  -2.80%  "
Set-TextValue $ws "E33" "  -0.57%  "
Set-TextValue $ws "E34" "  -2.04%  "
Set-TextValue $ws "D35" "1.381.03"
Set-TextValue $ws "E35" "  -0.92%  "
Set-TextValue $ws "E36" "  +2.75%  "
Set-TextValue $ws "E37" "  -2.69%  "
Set-TextValue $ws "E38" "  -0.26%  "
Set-TextValue $ws "D39" "2.65"
Set-TextValue $ws "E39" "  +2.11%  "
Set-TextValue $ws "D40" "0.0162"
Set-TextValue $ws "E40" "  -1.79%  "
Set-TextValue $ws "E41" "  -2.72%  "
Set-TextValue $ws "E42" "  +3.81%  "
Set-TextValue $ws "E43" "  -0.01%  "
Set-TextValue $ws "E44" "  -0.34%  "
Set-TextValue $ws "D45" "0.783"
Set-TextValue $ws "E45" "  -0.41%  "
Set-TextValue $ws "D46" "5.32"
Set-TextValue $ws "E46" "  -3.24%  "
Set-TextValue $ws "D47" "62.21"
Set-TextValue $ws "E47" "  -0.67%  "
Set-TextValue $ws "E48" "  -6.45%  "
Set-TextValue $ws "D49" "1.702.73"
Set-TextValue $ws "E49" "  +0.09%  "
Set-TextValue $ws "D50" "85.56"
Set-TextValue $ws "E50" "  -0.65%  "
Set-TextValue $ws "B51" "BabyDogeCoin"
Set-TextValue $ws "C51" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws "D51" "0.0₆0101"
Set-TextValue $ws "E51" "  -0.86%  "
